# fix: fixed title index storing in temp_book
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell E1 held the placeholder shared string "undefined";
# correct it to the proper Ukrainian column title.
$ws.Range("E1").Value = "Індекс посади"

# Move/restore the active selection cell on the sheet.
$ws.Range("B8").Select()
